# Auto-generated Excel COM-interop script applying the scheduled-runner update
# to the Gilgamesh_Profits workbook. Updates cached profit/price figures across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 195.77777
$ws.Range("I33").Value = 220.57143
$ws.Range("K33").Value = 220.57143
$ws.Range("M33").Value = 8.428570000000008
$ws.Range("H40").Value = 5123.1113
$ws.Range("J40").Value = 4935.5386
$ws.Range("L40").Value = 4935.5386
$ws.Range("N40").Value = -5285.5386
$ws.Range("H55").Value = 484.44446
$ws.Range("I55").Value = 790.3333
$ws.Range("J55").Value = 331.5
$ws.Range("K55").Value = 790.3333
$ws.Range("L55").Value = 331.5
$ws.Range("M55").Value = -576.3333
$ws.Range("N55").Value = -759.5
$ws.Range("H62").Value = 4123
$ws.Range("I62").Value = 4123
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4123
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3499
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4123
$ws.Range("I65").Value = 4123
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 20615
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -17495
$ws.Range("N65").Value = 0
$ws.Range("H88").Value = 7695961.5
$ws.Range("J88").Value = 5299.6
$ws.Range("L88").Value = 5299.6
$ws.Range("N88").Value = -6111.6
$ws.Range("H91").Value = 7695961.5
$ws.Range("J91").Value = 5299.6
$ws.Range("L91").Value = 5299.6
$ws.Range("N91").Value = -8107.6
$ws.Range("H92").Value = 2066.125
$ws.Range("H103").Value = 4763318.5
$ws.Range("I103").Value = 2058.889
$ws.Range("J103").Value = 11905208
$ws.Range("K103").Value = 6176.667
$ws.Range("L103").Value = 35715624
$ws.Range("M103").Value = -5590.667
$ws.Range("N103").Value = -35716796
$ws.Range("H112").Value = 865.1667
$ws.Range("J112").Value = 1595.5
$ws.Range("L112").Value = 4786.5
$ws.Range("N112").Value = -7002.5
$ws.Range("H116").Value = 4262.375
$ws.Range("I116").Value = 4245.5386
$ws.Range("K116").Value = 4245.5386
$ws.Range("M116").Value = -803.5385999999999
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 7329.143
$ws.Range("I132").Value = 7645.6
$ws.Range("K132").Value = 22936.8
$ws.Range("M132").Value = -20406.8
$ws.Range("H135").Value = 768.5
$ws.Range("I135").Value = 374.53333
$ws.Range("K135").Value = 3370.79997
$ws.Range("M135").Value = -835.79997
$ws.Range("H138").Value = 2767.9124
$ws.Range("I138").Value = 1082.3334
$ws.Range("J138").Value = 2966.2156
$ws.Range("K138").Value = 3247.0002
$ws.Range("L138").Value = 8898.6468
$ws.Range("M138").Value = 1892.9998
$ws.Range("N138").Value = -19178.6468

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 36245.25
$ws.Range("J43").Value = 36245.25
$ws.Range("L43").Value = 36245.25
$ws.Range("N43").Value = -36871.25
$ws.Range("H61").Value = 4369
$ws.Range("J61").Value = 4823
$ws.Range("L61").Value = 4823
$ws.Range("N61").Value = -5247
$ws.Range("H74").Value = 91183.08
$ws.Range("I74").Value = 104272.3
$ws.Range("J74").Value = 2830.875
$ws.Range("K74").Value = 104272.3
$ws.Range("L74").Value = 2830.875
$ws.Range("M74").Value = -103398.3
$ws.Range("N74").Value = -4578.875
$ws.Range("H77").Value = 91183.08
$ws.Range("I77").Value = 104272.3
$ws.Range("J77").Value = 2830.875
$ws.Range("K77").Value = 521361.5
$ws.Range("L77").Value = 14154.375
$ws.Range("M77").Value = -516993.5
$ws.Range("N77").Value = -22890.375
$ws.Range("H122").Value = 2613.7727
$ws.Range("I122").Value = 2523.9524
$ws.Range("K122").Value = 7571.8572
$ws.Range("M122").Value = -5121.8572
$ws.Range("H132").Value = 2259.5454
$ws.Range("I132").Value = 1446.3636
$ws.Range("K132").Value = 4339.0908
$ws.Range("M132").Value = -1809.0908
$ws.Range("H136").Value = 4369
$ws.Range("J136").Value = 4823
$ws.Range("L136").Value = 14469
$ws.Range("N136").Value = -19569

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 27781712
$ws.Range("I20").Value = 30868330
$ws.Range("J20").Value = 2150.6667
$ws.Range("K20").Value = 30868330
$ws.Range("L20").Value = 2150.6667
$ws.Range("M20").Value = -30868083
$ws.Range("N20").Value = -2644.6667
$ws.Range("H86").Value = 3856.2104
$ws.Range("I86").Value = 3574.5386
$ws.Range("J86").Value = 4466.5
$ws.Range("K86").Value = 3574.5386
$ws.Range("L86").Value = 4466.5
$ws.Range("M86").Value = -2451.5386
$ws.Range("N86").Value = -6712.5
$ws.Range("H89").Value = 3856.2104
$ws.Range("I89").Value = 3574.5386
$ws.Range("J89").Value = 4466.5
$ws.Range("K89").Value = 17872.693
$ws.Range("L89").Value = 22332.5
$ws.Range("M89").Value = -12256.693
$ws.Range("N89").Value = -33564.5
$ws.Range("H107").Value = 4274926.5
$ws.Range("I107").Value = 6411606
$ws.Range("J107").Value = 1567.8334
$ws.Range("K107").Value = 6411606
$ws.Range("L107").Value = 1567.8334
$ws.Range("M107").Value = -6409686
$ws.Range("N107").Value = -5407.8334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I22").Value = 298.8
$ws.Range("J22").Value = 299.5
$ws.Range("K22").Value = 298.8
$ws.Range("L22").Value = 299.5
$ws.Range("M22").Value = 51.19999999999999
$ws.Range("N22").Value = -999.5
$ws.Range("H31").Value = 4676.2324
$ws.Range("I31").Value = 3330.577
$ws.Range("J31").Value = 6734.294
$ws.Range("K31").Value = 3330.577
$ws.Range("L31").Value = 6734.294
$ws.Range("M31").Value = -3035.577
$ws.Range("N31").Value = -7324.294
$ws.Range("H34").Value = 4676.2324
$ws.Range("I34").Value = 3330.577
$ws.Range("J34").Value = 6734.294
$ws.Range("K34").Value = 3330.577
$ws.Range("L34").Value = 6734.294
$ws.Range("M34").Value = -3128.577
$ws.Range("N34").Value = -7138.294
$ws.Range("H132").Value = 15879714
$ws.Range("I132").Value = 6589.1
$ws.Range("J132").Value = 30309826
$ws.Range("K132").Value = 19767.3
$ws.Range("L132").Value = 90929478
$ws.Range("M132").Value = -17237.3
$ws.Range("N132").Value = -90934538
$ws.Range("H134").Value = 10748.333
$ws.Range("I134").Value = 9898
$ws.Range("K134").Value = 29694
$ws.Range("M134").Value = -27159

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 748.5454999999999
$ws.Range("I5").Value = 697.4286
$ws.Range("J5").Value = 838
$ws.Range("K5").Value = 2092.2858
$ws.Range("L5").Value = 2514
$ws.Range("M5").Value = -1980.2858
$ws.Range("N5").Value = -2738
$ws.Range("H54").Value = 4454.9
$ws.Range("I54").Value = 2933
$ws.Range("J54").Value = 5107.143
$ws.Range("K54").Value = 8799
$ws.Range("L54").Value = 15321.429
$ws.Range("M54").Value = -8240
$ws.Range("N54").Value = -16439.429
$ws.Range("H109").Value = 4172.4707
$ws.Range("I109").Value = 2186.4
$ws.Range("K109").Value = 6559.200000000001
$ws.Range("M109").Value = -5519.200000000001
$ws.Range("H112").Value = 2894.4
$ws.Range("I112").Value = 2014.5
$ws.Range("J112").Value = 3481
$ws.Range("K112").Value = 6043.5
$ws.Range("L112").Value = 10443
$ws.Range("M112").Value = -4935.5
$ws.Range("N112").Value = -12659
$ws.Range("H122").Value = 1694.6471
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1694.6471
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 15251.8239
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -20151.8239
$ws.Range("H131").Value = 10418872
$ws.Range("J131").Value = 2282.1333
$ws.Range("L131").Value = 6846.3999
$ws.Range("N131").Value = -16926.3999
$ws.Range("H135").Value = 748.5454999999999
$ws.Range("I135").Value = 697.4286
$ws.Range("J135").Value = 838
$ws.Range("K135").Value = 6276.8574
$ws.Range("L135").Value = 7542
$ws.Range("M135").Value = -3741.8574
$ws.Range("N135").Value = -12612

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 5595247
$ws.Range("J15").Value = 44652.625
$ws.Range("L15").Value = 44652.625
$ws.Range("N15").Value = -45228.625
$ws.Range("H31").Value = 1850
$ws.Range("I31").Value = 1850
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1850
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1558
$ws.Range("N31").ClearContents()
$ws.Range("H37").Value = 1850
$ws.Range("I37").Value = 1850
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1850
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -1573
$ws.Range("N37").ClearContents()
$ws.Range("H70").Value = 45460070
$ws.Range("I70").Value = 125003520
$ws.Range("J70").Value = 6670.5713
$ws.Range("K70").Value = 125003520
$ws.Range("L70").Value = 6670.5713
$ws.Range("M70").Value = -125003250
$ws.Range("N70").Value = -7210.5713
$ws.Range("H73").Value = 45460070
$ws.Range("I73").Value = 125003520
$ws.Range("J73").Value = 6670.5713
$ws.Range("K73").Value = 125003520
$ws.Range("L73").Value = 6670.5713
$ws.Range("M73").Value = -125002584
$ws.Range("N73").Value = -8542.5713
$ws.Range("H80").Value = 52634936
$ws.Range("I80").Value = 90912010
$ws.Range("J80").Value = 3967.375
$ws.Range("K80").Value = 90912010
$ws.Range("L80").Value = 3967.375
$ws.Range("M80").Value = -90911012
$ws.Range("N80").Value = -5963.375
$ws.Range("H81").Value = 5595247
$ws.Range("J81").Value = 44652.625
$ws.Range("L81").Value = 44652.625
$ws.Range("N81").Value = -46648.625
$ws.Range("H83").Value = 52634936
$ws.Range("I83").Value = 90912010
$ws.Range("J83").Value = 3967.375
$ws.Range("K83").Value = 454560050
$ws.Range("L83").Value = 19836.875
$ws.Range("M83").Value = -454555058
$ws.Range("N83").Value = -29820.875
$ws.Range("H84").Value = 5595247
$ws.Range("J84").Value = 44652.625
$ws.Range("L84").Value = 133957.875
$ws.Range("N84").Value = -143941.875
$ws.Range("H132").Value = 2453.5625
$ws.Range("I132").Value = 2205.182
$ws.Range("K132").Value = 6615.545999999999
$ws.Range("M132").Value = -4085.545999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1848
$ws.Range("I16").Value = 663.1667
$ws.Range("J16").Value = 3032.8333
$ws.Range("K16").Value = 663.1667
$ws.Range("L16").Value = 3032.8333
$ws.Range("M16").Value = -493.1667
$ws.Range("N16").Value = -3372.8333
$ws.Range("H46").Value = 3846.4285
$ws.Range("I46").Value = 3085
$ws.Range("J46").Value = 5750
$ws.Range("K46").Value = 3085
$ws.Range("L46").Value = 5750
$ws.Range("M46").Value = -2897
$ws.Range("N46").Value = -6126
$ws.Range("H132").Value = 8915.954
$ws.Range("J132").Value = 11348
$ws.Range("L132").Value = 34044
$ws.Range("N132").Value = -39104
$ws.Range("H136").Value = 5156.421
$ws.Range("I136").Value = 4863.364
$ws.Range("J136").Value = 5559.375
$ws.Range("K136").Value = 14590.092
$ws.Range("L136").Value = 16678.125
$ws.Range("M136").Value = -12040.092
$ws.Range("N136").Value = -21778.125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5961.7334
$ws.Range("I81").Value = 5928.375
$ws.Range("K81").Value = 11856.75
$ws.Range("M81").Value = -10795.75
$ws.Range("H84").Value = 5961.7334
$ws.Range("I84").Value = 5928.375
$ws.Range("K84").Value = 59283.75
$ws.Range("M84").Value = -53979.75
$ws.Range("H122").Value = 8067628.5
$ws.Range("I122").Value = 3295.5925
$ws.Range("J122").Value = 62501876
$ws.Range("K122").Value = 9886.7775
$ws.Range("L122").Value = 187505628
$ws.Range("M122").Value = -7436.7775
$ws.Range("N122").Value = -187510528
$ws.Range("H132").Value = 2483.1304
$ws.Range("I132").Value = 2265.7
$ws.Range("J132").Value = 3932.6667
$ws.Range("K132").Value = 6797.099999999999
$ws.Range("L132").Value = 11798.0001
$ws.Range("M132").Value = -4267.099999999999
$ws.Range("N132").Value = -16858.0001

